# Update "想去人数" (number of people interested) values (column F)
# on the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 744
$ws1.Range("F4").Value = 254
$ws1.Range("F5").Value = 3173
$ws1.Range("F6").Value = 64
$ws1.Range("F7").Value = 3919
$ws1.Range("F8").Value = 481
$ws1.Range("F9").Value = 982

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 744
$ws4.Range("F5").Value = 254
$ws4.Range("F6").Value = 3173
$ws4.Range("F7").Value = 64
$ws4.Range("F8").Value = 3919
$ws4.Range("F9").Value = 481
$ws4.Range("F10").Value = 982
